$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-RowData {
    param($ws, $rowA, $rowB, $cols)
    $valsA = @{}
    $valsB = @{}
    foreach ($col in $cols) {
        $valsA[$col] = $ws.Range($col + $rowA).Value()
        $valsB[$col] = $ws.Range($col + $rowB).Value()
    }
    foreach ($col in $cols) {
        $va = $valsA[$col]
        $vb = $valsB[$col]
        if ($vb -ne $null) { $ws.Range($col + $rowA).Value = $vb } else { $ws.Range($col + $rowA).ClearContents() }
        if ($va -ne $null) { $ws.Range($col + $rowB).Value = $va } else { $ws.Range($col + $rowB).ClearContents() }
    }
}

# --- 1) Two pairs of rows were re-sorted (same date, order of the two matches swapped) ---
Swap-RowData $ws 83 85 $dataCols
Swap-RowData $ws 86 87 $dataCols

# --- 2) Two new finished matches are inserted before the still-not-played id=7809824 match, ---
#        which shifts from row 205 down to row 207.

# 2a) Give rows 205 and 206 the same per-column formatting as the existing data rows
#     (column A bold/center/bordered style, column E date style) *before* writing any
#     values into them, reusing the existing style indices (no new styles created).
$ws.Range("A204").Copy()
$ws.Range("A206").PasteSpecial(-4122)
$ws.Range("E204").Copy()
$ws.Range("E206").PasteSpecial(-4122)

# 2b) Move the old row 205 (id 7809824, still unplayed) down to row 207, keeping its
#     formatting, but with the running index (column A) updated from 203 to 205.
$ws.Range("A205").Copy()
$ws.Range("A207").PasteSpecial(-4122)
$ws.Range("E205").Copy()
$ws.Range("E207").PasteSpecial(-4122)

foreach ($col in $dataCols) {
    $v = $ws.Range($col + "205").Value()
    if ($v -ne $null) {
        $ws.Range($col + "207").Value = $v
    }
}
$ws.Range("A207").Value = 205

# 2c) Clear the old row 205 contents (formatting stays) so it can receive the new match.
$ws.Range("A205:AC205").ClearContents()

# 2d) Write the two new, already-finished matches into rows 205 and 206.
$ws.Range("A205").Value = 204
$ws.Range("B205").Value = 7011608
$ws.Range("C205").Value = "Azerbaijan Premier League"
$ws.Range("D205").Value = "Azerbaijan Premier League"
$ws.Range("E205").Value = 45346.33333333334
$ws.Range("F205").Value = "PFK Turan Tovuz"
$ws.Range("G205").Value = "Neftchi Baku"
$ws.Range("H205").Value = 1
$ws.Range("I205").Value = 1
$ws.Range("J205").Value = "D"
$ws.Range("K205").Value = 3.4
$ws.Range("L205").Value = 3.4
$ws.Range("M205").Value = 1.909
$ws.Range("N205").Value = 2.9
$ws.Range("O205").Value = 3.3
$ws.Range("P205").Value = 2.25
$ws.Range("Q205").Value = 0.25
$ws.Range("R205").Value = 1.8
$ws.Range("S205").Value = 2
$ws.Range("T205").Value = 2.25
$ws.Range("U205").Value = 1.975
$ws.Range("V205").Value = 1.825
$ws.Range("W205").Value = -1
$ws.Range("X205").Value = 2.3
$ws.Range("Y205").Value = -1
$ws.Range("Z205").Value = 0.4
$ws.Range("AA205").Value = -0.5
$ws.Range("AB205").Value = -0.5
$ws.Range("AC205").Value = 0.4125

$ws.Range("A206").Value = 205
$ws.Range("B206").Value = 7011611
$ws.Range("C206").Value = "Azerbaijan Premier League"
$ws.Range("D206").Value = "Azerbaijan Premier League"
$ws.Range("E206").Value = 45346.4375
$ws.Range("F206").Value = "Sabail FC"
$ws.Range("G206").Value = "Sabah"
$ws.Range("H206").Value = 2
$ws.Range("I206").Value = 0
$ws.Range("J206").Value = "H"
$ws.Range("K206").Value = 2.75
$ws.Range("L206").Value = 3.2
$ws.Range("M206").Value = 2.3
$ws.Range("N206").Value = 3.3
$ws.Range("O206").Value = 3.2
$ws.Range("P206").Value = 2
$ws.Range("Q206").Value = 0.25
$ws.Range("R206").Value = 2
$ws.Range("S206").Value = 1.8
$ws.Range("T206").Value = 2.75
$ws.Range("U206").Value = 1.925
$ws.Range("V206").Value = 1.775
$ws.Range("W206").Value = 2.3
$ws.Range("X206").Value = -1
$ws.Range("Y206").Value = -1
$ws.Range("Z206").Value = 1
$ws.Range("AA206").Value = -1
$ws.Range("AB206").Value = -1
$ws.Range("AC206").Value = 0.7749999999999999

"done"
